# Applies the "想去人数" (F column) count updates captured in the commit
# "Update gh-pages to output generated at 456a3b4" across all four sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1180
$ws.Range("F5").Value = 9271
$ws.Range("F8").Value = 7206
$ws.Range("F9").Value = 192
$ws.Range("F10").Value = 326
$ws.Range("F11").Value = 5655
$ws.Range("F12").Value = 78
$ws.Range("F13").Value = 6491
$ws.Range("F14").Value = 1107
$ws.Range("F15").Value = 445
$ws.Range("F18").Value = 337
$ws.Range("F19").Value = 285
$ws.Range("F22").Value = 163
$ws.Range("F23").Value = 10487
$ws.Range("F26").Value = 1997
$ws.Range("F27").Value = 2349
$ws.Range("F28").Value = 48
$ws.Range("F29").Value = 2264
$ws.Range("F30").Value = 86
$ws.Range("F31").Value = 88
$ws.Range("F32").Value = 195
$ws.Range("F33").Value = 26
$ws.Range("F34").Value = 2160
$ws.Range("F35").Value = 316
$ws.Range("F36").Value = 1444
$ws.Range("F38").Value = 5384
$ws.Range("F39").Value = 446
$ws.Range("F40").Value = 1216
$ws.Range("F41").Value = 723
$ws.Range("F42").Value = 129
$ws.Range("F44").Value = 1109
$ws.Range("F45").Value = 1083
$ws.Range("F46").Value = 1009
$ws.Range("F47").Value = 1410
$ws.Range("F48").Value = 70
$ws.Range("F49").Value = 1106

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 19
$ws.Range("F6").Value = 3
$ws.Range("F10").Value = 197
$ws.Range("F11").Value = 1
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 100
$ws.Range("F16").Value = 1
$ws.Range("F18").Value = 913
$ws.Range("F20").Value = 7
$ws.Range("F22").Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 74

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1180
$ws.Range("F5").Value = 9271
$ws.Range("F6").Value = 7206
$ws.Range("F7").Value = 192
$ws.Range("F8").Value = 74
$ws.Range("F10").Value = 3
$ws.Range("F11").Value = 5655
$ws.Range("F12").Value = 5655
$ws.Range("F13").Value = 78
$ws.Range("F14").Value = 6491
$ws.Range("F15").Value = 6491
$ws.Range("F17").Value = 445
$ws.Range("F19").Value = 623
$ws.Range("F20").Value = 337
$ws.Range("F21").Value = 285
$ws.Range("F24").Value = 241
$ws.Range("F25").Value = 163
$ws.Range("F26").Value = 197
$ws.Range("F27").Value = 10487
$ws.Range("F28").Value = 1997
$ws.Range("F29").Value = 2349
$ws.Range("F30").Value = 48
$ws.Range("F31").Value = 2264
$ws.Range("F34").Value = 195
$ws.Range("F36").Value = 2160
$ws.Range("F37").Value = 316
$ws.Range("F38").Value = 1444
$ws.Range("F40").Value = 5384
$ws.Range("F41").Value = 446
$ws.Range("F42").Value = 1216
$ws.Range("F43").Value = 723
$ws.Range("F45").Value = 169
$ws.Range("F47").Value = 1083
$ws.Range("F48").Value = 1009
$ws.Range("F49").Value = 1410
$ws.Range("F50").Value = 70
$ws.Range("F51").Value = 1106
